$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.878.31'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '2.955.99'
$ws.Range("E3").Value = '  +2.35%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'352.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'111.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").Value = "'0.565"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'0.634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("D10").Value = "'39.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.34%  '
$ws.Range("D11").Value = "'0.0896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.06%  '
$ws.Range("D12").Value = "'0.136"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").Value = "'19.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = "'8.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").Value = '3.425.08'
$ws.Range("D16").Value = '2.961.32'
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("D17").Value = "'0.999"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '52.014.52'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("D20").Value = "'14.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.67%  '
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("D22").Value = '0.0₃0991'
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("D23").Value = "'71.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("D24").Value = "'272.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").Value = "'2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = "'0.181"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.32%  '
$ws.Range("D27").Value = "'27.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = "'7.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +18.40%  '
$ws.Range("E30").Value = '  +22.69%  '
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = "'6.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.62%  '
$ws.Range("D33").Value = "'37.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.32%  '
$ws.Range("D34").Value = "'53.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = "'3.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("D38").Value = "'1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -18.78%  '
$ws.Range("D39").Value = "'18.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").Value = "'2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("E42").Value = '  +2.21%  '
$ws.Range("D43").Value = "'23.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("E46").Value = '  +1.93%  '
$ws.Range("D47").Value = '2.167.55'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").Value = "'114.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.87%  '
$ws.Range("D49").Value = "'0.245"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.98%  '
$ws.Range("D50").Value = "'0.0340"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.14%  '
$ws.Range("D51").Value = "'0.933"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.76%  '
